$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text so values such as "1.007" or
# "8.000" are not reinterpreted as numbers (which would drop precision and
# trailing zeros). This mirrors the original file, where every Price/Volume
# cell is an inline/shared string.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.462.83'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '1.998.87'
$ws.Range('E3').Value = '  -5.73%  '
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '331.29'
$ws.Range('E5').Value = '  -4.42%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.4944'
$ws.Range('E7').Value = '  -4.88%  '
$ws.Range('D8').Value = '0.4187'
$ws.Range('E8').Value = '  -6.21%  '
$ws.Range('D9').Value = '53.34'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '0.08822'
$ws.Range('E10').Value = '  -5.98%  '
$ws.Range('D11').Value = '1.115'
$ws.Range('E11').Value = '  -5.86%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '2.049.49'
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '23.15'
$ws.Range('E13').Value = '  -8.83%  '
$ws.Range('D14').Value = '8.000'
$ws.Range('E14').Value = '  -6.79%  '
$ws.Range('D15').Value = '6.476'
$ws.Range('E15').Value = '  -6.90%  '
$ws.Range('D16').Value = '96.27'
$ws.Range('E16').Value = '  -6.49%  '
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '0.00001108'
$ws.Range('E18').Value = '  -4.84%  '
$ws.Range('D19').Value = '0.06630'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('D20').Value = '19.62'
$ws.Range('E20').Value = '  -8.96%  '
$ws.Range('D21').Value = '1.008'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '5.969'
$ws.Range('E22').Value = '  -5.24%  '
$ws.Range('D23').Value = '29.497.89'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').Value = '11.83'
$ws.Range('E24').Value = '  -7.07%  '
$ws.Range('D25').Value = '2.285'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').Value = '2.307.05'
$ws.Range('E26').Value = '  -2.15%  '
$ws.Range('D27').Value = '6.673'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').Value = '157.28'
$ws.Range('E28').Value = '  -3.28%  '
$ws.Range('D29').Value = '20.52'
$ws.Range('E29').Value = '  -7.23%  '
$ws.Range('D30').Value = '2.351'
$ws.Range('E30').Value = '  -7.54%  '
$ws.Range('D31').Value = '126.85'
$ws.Range('E31').Value = '  -5.35%  '
$ws.Range('D32').Value = '1.053'
$ws.Range('E32').Value = '  -8.82%  '
$ws.Range('D33').Value = '0.09919'
$ws.Range('E33').Value = '  -6.08%  '
$ws.Range('D34').Value = '1.555'
$ws.Range('E34').Value = '  -13.27%  '
$ws.Range('D35').Value = '5.832'
$ws.Range('E35').Value = '  -6.98%  '
$ws.Range('D36').Value = '3.771'
$ws.Range('E36').Value = '  -4.92%  '
$ws.Range('D37').Value = '9.574'
$ws.Range('E37').Value = '  -11.34%  '
$ws.Range('D38').Value = '0.02449'
$ws.Range('E38').Value = '  -6.45%  '
$ws.Range('D39').Value = '0.06387'
$ws.Range('E39').Value = '  -7.01%  '
$ws.Range('D40').Value = '1.285'
$ws.Range('E40').Value = '  -3.80%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '11.74'
$ws.Range('E41').Value = '  -7.74%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6491'
$ws.Range('E42').Value = '  -8.50%  '
$ws.Range('D43').Value = '0.2066'
$ws.Range('E43').Value = '  -8.03%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = '0.6324'
$ws.Range('E45').Value = '  -7.70%  '
$ws.Range('D46').Value = '13.37'
$ws.Range('E46').Value = '  -9.06%  '
$ws.Range('D47').Value = '2.199'
$ws.Range('E47').Value = '  -7.36%  '
$ws.Range('D48').Value = '1.272'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = '3.545'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = '0.00000000335'
$ws.Range('E50').Value = '  -5.90%  '
$ws.Range('D51').Value = '0.07005'
$ws.Range('E51').Value = '  -1.66%  '
